# Update cryptos price/volume data per Thu Mar 14 13:11:26 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '72.801.61'
$ws.Cells.Item(2, 5).Value = '  -0.09%  '
$ws.Cells.Item(3, 4).Value = '3.944.07'
$ws.Cells.Item(3, 5).Value = '  -1.76%  '
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '604.77'
$c.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +1.97%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '171.82'
$c.ClearFormats()
$ws.Cells.Item(6, 5).Value = '  +12.56%  '
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.685'
$c.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  -0.40%  '
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Cells.Item(8, 5).Value = '  +0.08%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.787'
$c.ClearFormats()
$ws.Cells.Item(9, 5).Value = '  +3.50%  '
$ws.Cells.Item(10, 5).Value = '  +9.32%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '56.13'
$c.ClearFormats()
$ws.Cells.Item(11, 5).Value = '  +3.15%  '
$ws.Cells.Item(12, 5).Value = '  +3.52%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '11.57'
$c.ClearFormats()
$ws.Cells.Item(13, 5).Value = '  +5.30%  '
$ws.Cells.Item(14, 4).Value = '4.572.48'
$ws.Cells.Item(14, 5).Value = '  -1.85%  '
$ws.Cells.Item(15, 2).Value = 'WrappedEther'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(15, 4).Value = '3.948.24'
$ws.Cells.Item(15, 5).Value = '  -1.99%  '
$ws.Cells.Item(16, 2).Value = 'Chainlink'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '21.40'
$c.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  +3.62%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '14.23'
$c.ClearFormats()
$ws.Cells.Item(17, 5).Value = '  -0.39%  '
$ws.Cells.Item(18, 5).Value = '  -3.35%  '
$ws.Cells.Item(19, 4).Value = '72.713.39'
$ws.Cells.Item(19, 5).Value = '  +0.03%  '
$ws.Cells.Item(20, 5).Value = '  -1.18%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '449.53'
$c.ClearFormats()
$ws.Cells.Item(21, 5).Value = '  +2.53%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '4.82'
$c.ClearFormats()
$ws.Cells.Item(22, 5).Value = '  +0.61%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '95.74'
$c.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  -1.50%  '
$ws.Cells.Item(24, 5).Value = '  -5.46%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '14.09'
$c.ClearFormats()
$ws.Cells.Item(25, 5).Value = '  -1.89%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '4.30'
$c.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -1.18%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '11.20'
$c.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  -2.53%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '5.96'
$c.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  +0.65%  '
$ws.Cells.Item(29, 5).Value = '  -3.38%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '35.88'
$c.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  -2.15%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '7.84'
$c.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -1.07%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '13.88'
$c.ClearFormats()
$ws.Cells.Item(32, 5).Value = '  +1.77%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '50.71'
$c.ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +1.13%  '
$ws.Cells.Item(34, 5).Value = '  -4.28%  '
$ws.Cells.Item(35, 5).Value = '  +14.76%  '
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '68.95'
$c.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  -3.04%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '632.84'
$c.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -7.75%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '0.428'
$c.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -4.09%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '3.37'
$c.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +0.54%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  +0.17%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '0.146'
$c.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  -2.01%  '
$ws.Cells.Item(42, 2).Value = 'dogwifhat'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '3.30'
$c.ClearFormats()
$ws.Cells.Item(42, 5).Value = '  +43.97%  '
$ws.Cells.Item(43, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Cells.Item(43, 5).Value = '  +0.01%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.0479'
$c.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  -2.97%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '10.57'
$c.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -5.56%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '0.148'
$c.ClearFormats()
$ws.Cells.Item(46, 5).Value = '  -2.18%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '2.63'
$c.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  -4.86%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '3.38'
$c.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  -0.36%  '
$ws.Cells.Item(49, 5).Value = '  -16.29%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.000281'
$c.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  +4.07%  '
$ws.Cells.Item(51, 4).Value = '2.833.46'
$ws.Cells.Item(51, 5).Value = '  +0.23%  '
